$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Andrea Favero
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = "-"
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 15
$ws.Range("G2").Value = 15

# Row 3 - Eleonora Thiella
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = 20
$ws.Range("F3").Value = 15
$ws.Range("G3").Value = 3

# Row 4 - Federico Caldart
$ws.Range("B4").Value = "-"
$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = 20
$ws.Range("F4").Value = 16
$ws.Range("G4").Value = 10

# Row 5 - Giovanni Cavallin
$ws.Range("B5").Value = "-"
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = "-"
$ws.Range("E5").Value = 15
$ws.Range("F5").Value = 14
$ws.Range("G5").Value = 16

# Row 6 - Giovanni Dalla Riva
$ws.Range("B6").Value = "-"
$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = 20
$ws.Range("F6").Value = 10
$ws.Range("G6").Value = 16

# Row 7 - Lorenzo Menegon
$ws.Range("B7").Value = 5
$ws.Range("C7").Value = "-"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = 15
$ws.Range("G7").Value = 10

# Row 8 - Stefano Panozzo
$ws.Range("B8").Value = "-"
$ws.Range("C8").Value = "-"
$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = 20
$ws.Range("F8").Value = 15
$ws.Range("G8").Value = 10

# Update chart position/size (moved up-left and enlarged)
$chart = $ws.ChartObjects().Item(1)
$chart.Left = 689.9648043799212
$chart.Top = 0.0
$chart.Width = 737.8124409448818
$chart.Height = 286.35

# Update selection
$ws.Range("A1:H9").Select()
